# "Added Unittests and documentation"
#
# Fills in the unit-test-plan worksheet for the LibraryItem class:
#   - records the developer's name
#   - fills in Preconditions / Method Inputs / Expected Result for the
#     __init__ exception test cases (rows 7-10)
#   - fills in Preconditions / Method Inputs / Expected Result for the
#     attribute-getter test cases (rows 11-13)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header: developer name -------------------------------------------------
$ws.Range("C3").Value = 'Navkaran Singh Sidhu'

# --- __init__ exception test cases (rows 7-10) ------------------------------
# Row 7: title blank -> ValueError-style assertion on Argument Values
$ws.Range("E7").Value = 'None'
$ws.Range("F7").Value = '"Kingdom of Ash"                                                             "al sweigart "                                                                                  Genre.TRUE_CRIME'
$ws.Range("G7").Value = 'Attribute set to Argument Values.'

# Row 8: author blank -> ValueError
$ws.Range("E8").Value = 'None'
$ws.Range("F8").Value = '""                                                                                                           "al sweigart"                                                                                    Genre.TRUE_CRIME'
$ws.Range("G8").Value = 'ValueError'

# Row 9: title blank (second variant) -> ValueError
$ws.Range("E9").Value = 'None'
$ws.Range("F9").Value = '"Kingdom of Ash"                                                                     ""                                                                                                            Genre.TRUE_CRIME'
$ws.Range("G9").Value = 'ValueError'

# Row 10: invalid Genre -> ValueError
$ws.Range("E10").Value = 'None'
$ws.Range("F10").Value = '"Kingdom of Ash"                                                           "al sweigart "                                                                                  "INVALID"'
$ws.Range("G10").Value = 'ValueError'

# --- Attribute getter test cases (rows 11-13) -------------------------------
# Row 11: returns title attribute
$ws.Range("E11").Value = 'LibraryItem("Kingdom of Ash"                                                                 "al sweigart "                                                                                  Genre.TRUE_CRIME)'
$ws.Range("F11").Value = 'None'
$ws.Range("G11").Value = '"Kingdom of Ash "'

# Row 12: returns author attribute
$ws.Range("E12").Value = 'LibraryItem("Kingdom of Ash"                                          "al sweigart "                                                                                  Genre.TRUE_CRIME)'
$ws.Range("F12").Value = 'None'
$ws.Range("G12").Value = '"al sweigart"'

# Row 13: returns Genre attribute
$ws.Range("E13").Value = 'LibraryItem("Kingdom of Ash"                                          "al sweigart "                                                                                  Genre.TRUE_CRIME)'
$ws.Range("F13").Value = 'None'
$ws.Range("G13").Value = 'Genre.TRUE_CRIME'

# --- Leave the cursor where the author left it (cell E13) ------------------
$ws.Range("E13").Select()
